$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3980
$ws.Range("J7").Value = 3980
$ws.Range("L7").Value = 3980
$ws.Range("N7").Value = -4204
$ws.Range("H14").Value = 3980
$ws.Range("J14").Value = 3980
$ws.Range("L14").Value = 3980
$ws.Range("N14").Value = -4362
$ws.Range("H19").Value = 1923.3529
$ws.Range("J19").Value = 489.83334
$ws.Range("L19").Value = 489.83334
$ws.Range("N19").Value = -839.83334
$ws.Range("H40").Value = 1160.8636
$ws.Range("I40").Value = 799.8
$ws.Range("J40").Value = 1267.0588
$ws.Range("K40").Value = 799.8
$ws.Range("L40").Value = 1267.0588
$ws.Range("M40").Value = -624.8
$ws.Range("N40").Value = -1617.0588
$ws.Range("H62").Value = 8332.75
$ws.Range("I62").Value = 10413.286
$ws.Range("J62").Value = 5420
$ws.Range("K62").Value = 10413.286
$ws.Range("L62").Value = 5420
$ws.Range("M62").Value = -9789.286
$ws.Range("N62").Value = -6668
$ws.Range("H65").Value = 8332.75
$ws.Range("I65").Value = 10413.286
$ws.Range("J65").Value = 5420
$ws.Range("K65").Value = 52066.43
$ws.Range("L65").Value = 27100
$ws.Range("M65").Value = -48946.43
$ws.Range("N65").Value = -33340
$ws.Range("H111").Value = 962.1
$ws.Range("I111").Value = 727.1429000000001
$ws.Range("J111").Value = 1510.3334
$ws.Range("K111").Value = 2181.4287
$ws.Range("L111").Value = 4531.0002
$ws.Range("M111").Value = 885.5712999999996
$ws.Range("N111").Value = -10665.0002

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 630
$ws.Range("I2").Value = 537.5
$ws.Range("K2").Value = 537.5
$ws.Range("M2").Value = -424.5
$ws.Range("H11").Value = 7628.5
$ws.Range("J11").Value = 7628.5
$ws.Range("L11").Value = 7628.5
$ws.Range("N11").Value = -7916.5
$ws.Range("H32").Value = 8248.724
$ws.Range("I32").Value = 4763.9766
$ws.Range("K32").Value = 4763.9766
$ws.Range("M32").Value = -4476.9766
$ws.Range("H45").Value = 1041.7778
$ws.Range("I45").Value = 875.1579
$ws.Range("J45").Value = 1437.5
$ws.Range("K45").Value = 875.1579
$ws.Range("L45").Value = 1437.5
$ws.Range("M45").Value = -498.1579
$ws.Range("N45").Value = -2191.5
$ws.Range("H61").Value = 2532.625
$ws.Range("I61").Value = 1688.6938
$ws.Range("J61").Value = 5289.467
$ws.Range("K61").Value = 1688.6938
$ws.Range("L61").Value = 5289.467
$ws.Range("M61").Value = -1476.6938
$ws.Range("N61").Value = -5713.467
$ws.Range("H74").Value = 1494.3802
$ws.Range("I74").Value = 978.678
$ws.Range("K74").Value = 978.678
$ws.Range("M74").Value = -104.678
$ws.Range("H77").Value = 1494.3802
$ws.Range("I77").Value = 978.678
$ws.Range("K77").Value = 4893.39
$ws.Range("M77").Value = -525.3900000000003
$ws.Range("H116").Value = 630
$ws.Range("I116").Value = 537.5
$ws.Range("K116").Value = 537.5
$ws.Range("M116").Value = 1756.5
$ws.Range("H136").Value = 2532.625
$ws.Range("I136").Value = 1688.6938
$ws.Range("J136").Value = 5289.467
$ws.Range("K136").Value = 5066.0814
$ws.Range("L136").Value = 15868.401
$ws.Range("M136").Value = -2516.0814
$ws.Range("N136").Value = -20968.401

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 630
$ws.Range("I3").Value = 537.5
$ws.Range("K3").Value = 537.5
$ws.Range("M3").Value = -423.5
$ws.Range("H99").Value = 5772131
$ws.Range("I99").Value = 2027776
$ws.Range("K99").Value = 2027776
$ws.Range("M99").Value = -2026278

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 722.619
$ws.Range("I22").Value = 412.08334
$ws.Range("J22").Value = 1136.6666
$ws.Range("K22").Value = 412.08334
$ws.Range("L22").Value = 1136.6666
$ws.Range("M22").Value = -62.08334000000002
$ws.Range("N22").Value = -1836.6666
$ws.Range("H31").Value = 2170.8877
$ws.Range("I31").Value = 1479.8524
$ws.Range("J31").Value = 3310.162
$ws.Range("K31").Value = 1479.8524
$ws.Range("L31").Value = 3310.162
$ws.Range("M31").Value = -1184.8524
$ws.Range("N31").Value = -3900.162
$ws.Range("H34").Value = 2170.8877
$ws.Range("I34").Value = 1479.8524
$ws.Range("J34").Value = 3310.162
$ws.Range("K34").Value = 1479.8524
$ws.Range("L34").Value = 3310.162
$ws.Range("M34").Value = -1277.8524
$ws.Range("N34").Value = -3714.162
$ws.Range("H107").Value = 277.08694
$ws.Range("I107").Value = 193.75
$ws.Range("K107").Value = 193.75
$ws.Range("M107").Value = 1726.25
$ws.Range("H134").Value = 2174.348
$ws.Range("I134").Value = 1345.5555
$ws.Range("J134").Value = 2707.1428
$ws.Range("K134").Value = 4036.6665
$ws.Range("L134").Value = 8121.428400000001
$ws.Range("M134").Value = -1501.6665
$ws.Range("N134").Value = -13191.4284

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 914.2857
$ws.Range("I92").Value = 900
$ws.Range("J92").Value = 933.3333
$ws.Range("K92").Value = 2700
$ws.Range("L92").Value = 2799.9999
$ws.Range("M92").Value = -1452
$ws.Range("N92").Value = -5295.9999

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3041.1
$ws.Range("J126").Value = 1933.3334
$ws.Range("L126").Value = 5800.0002
$ws.Range("N126").Value = -10740.0002
$ws.Range("H132").Value = 3478.0217
$ws.Range("I132").Value = 3302.7576
$ws.Range("J132").Value = 3922.923
$ws.Range("K132").Value = 9908.272799999999
$ws.Range("L132").Value = 11768.769
$ws.Range("M132").Value = -7378.272799999999
$ws.Range("N132").Value = -16828.769
$ws.Range("H136").Value = 6269.3716
$ws.Range("J136").Value = 6269.3716
$ws.Range("L136").Value = 18808.1148
$ws.Range("N136").Value = -23908.1148

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 446.58334
$ws.Range("I22").Value = 361.375
$ws.Range("K22").Value = 361.375
$ws.Range("M22").Value = -66.375
$ws.Range("H27").Value = 446.58334
$ws.Range("I27").Value = 361.375
$ws.Range("K27").Value = 361.375
$ws.Range("M27").Value = -254.375

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 41565.84
$ws.Range("I122").Value = 51448.3
$ws.Range("J122").Value = 2036
$ws.Range("K122").Value = 154344.9
$ws.Range("L122").Value = 6108
$ws.Range("M122").Value = -151894.9
$ws.Range("N122").Value = -11008
$ws.Range("H126").Value = 1395.375
$ws.Range("I126").Value = 1321.2
$ws.Range("J126").Value = 1519
$ws.Range("K126").Value = 3963.6
$ws.Range("L126").Value = 4557
$ws.Range("M126").Value = -1493.6
$ws.Range("N126").Value = -9497
$ws.Range("H132").Value = 19420.352
$ws.Range("I132").Value = 28797.861
$ws.Range("K132").Value = 86393.583
$ws.Range("M132").Value = -83863.583
$ws.Range("H136").Value = 1587.0492
$ws.Range("I136").Value = 1173.4242
$ws.Range("K136").Value = 3520.2726
$ws.Range("M136").Value = -970.2725999999998

Write-Output "edits applied"
